$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E2" = 3
    "F2" = 1
    "G2" = 170.8510486666667
    "H2" = 512.553146
    "I2" = 0.8284864843875368
    "J2" = 0.8284864843875368
    "M2" = 22.56839366666667
    "N2" = 67.705181
    "O2" = 0.9671969057210177
    "P2" = 0.9671969057210175
    "Q2" = 3855.833724672158
    "R2" = 34702.50352204942
    "S2" = 0.8013095641313098
    "T2" = 0.8013095641313096
    "E3" = 3
    "F3" = 1
    "G3" = 170.8510486666667
    "H3" = 512.553146
    "I3" = 0.8284864843875368
    "J3" = 0.8284864843875368
    "M3" = 0.2726246666666667
    "N3" = 0.817874
    "O3" = 0.01168367310131955
    "P3" = 0.01168367310131955
    "Q3" = 46.57821019240044
    "R3" = 419.203891731604
    "S3" = 0.009679765252445465
    "T3" = 0.009679765252445464
    "E4" = 3
    "F4" = 1
    "G4" = 170.8510486666667
    "H4" = 512.553146
    "I4" = 0.8284864843875368
    "J4" = 0.8284864843875368
    "K4" = 3
    "L4" = 1
    "M4" = 0.4927966666666667
    "N4" = 1.47839
    "O4" = 0.02111942117766284
    "P4" = 0.02111942117766283
    "Q4" = 84.19482727943777
    "R4" = 757.75344551494
    "S4" = 0.01749715500378158
    "T4" = 0.01749715500378157
    "G5" = 8.580369333333332
    "I5" = 0.04160770495233658
    "J5" = 0.04160770495233659
    "M5" = 22.56839366666667
    "N5" = 67.705181
    "O5" = 0.9671969057210177
    "P5" = 0.9671969057210175
    "Q5" = 193.6451529200608
    "R5" = 1742.806376280548
    "S5" = 0.04024284348405301
    "T5" = 0.040242843484053
    "G6" = 8.580369333333332
    "I6" = 0.04160770495233658
    "J6" = 0.04160770495233659
    "M6" = 0.2726246666666667
    "N6" = 0.817874
    "O6" = 0.01168367310131955
    "P6" = 0.01168367310131955
    "Q6" = 2.339220329376889
    "R6" = 21.052982964392
    "S6" = 0.0004861308231592553
    "T6" = 0.0004861308231592552
    "G7" = 8.580369333333332
    "I7" = 0.04160770495233658
    "J7" = 0.04160770495233659
    "K7" = 3
    "L7" = 1
    "M7" = 0.4927966666666667
    "N7" = 1.47839
    "O7" = 0.02111942117766284
    "P7" = 0.02111942117766283
    "Q7" = 4.228377406235555
    "R7" = 38.05539665612
    "S7" = 0.000878730645124324
    "T7" = 0.000878730645124324
    "G8" = 24.73238366666666
    "H8" = 74.19715099999999
    "I8" = 0.1199316349207643
    "J8" = 0.1199316349207643
    "M8" = 22.56839366666667
    "N8" = 67.705181
    "O8" = 0.9671969057210177
    "P8" = 0.9671969057210175
    "Q8" = 558.17017090437
    "R8" = 5023.53153813933
    "S8" = 0.115997506193426
    "T8" = 0.115997506193426
    "G9" = 24.73238366666666
    "H9" = 74.19715099999999
    "I9" = 0.1199316349207643
    "J9" = 0.1199316349207643
    "M9" = 0.2726246666666667
    "N9" = 0.817874
    "O9" = 0.01168367310131955
    "P9" = 0.01168367310131955
    "Q9" = 6.742657852997111
    "R9" = 60.68392067697399
    "S9" = 0.001401242016921011
    "T9" = 0.00140124201692101
    "G10" = 24.73238366666666
    "H10" = 74.19715099999999
    "I10" = 0.1199316349207643
    "J10" = 0.1199316349207643
    "K10" = 3
    "L10" = 1
    "M10" = 0.4927966666666667
    "N10" = 1.47839
    "O10" = 0.02111942117766284
    "P10" = 0.02111942117766283
    "Q10" = 12.18803622965444
    "R10" = 109.69232606689
    "S10" = 0.002532886710417317
    "T10" = 0.002532886710417317
    "G11" = 2.056881333333333
    "H11" = 6.170643999999999
    "I11" = 0.009974175739362347
    "J11" = 0.009974175739362347
    "M11" = 22.56839366666667
    "N11" = 67.705181
    "O11" = 0.9671969057210177
    "P11" = 0.9671969057210175
    "Q11" = 46.42050765628488
    "R11" = 417.7845689065639
    "S11" = 0.009646991912228906
    "T11" = 0.009646991912228904
    "G12" = 2.056881333333333
    "H12" = 6.170643999999999
    "I12" = 0.009974175739362347
    "J12" = 0.009974175739362347
    "M12" = 0.2726246666666667
    "N12" = 0.817874
    "O12" = 0.01168367310131955
    "P12" = 0.01168367310131955
    "Q12" = 0.5607565878728888
    "R12" = 5.046809290855999
    "S12" = 0.0001165350087938219
    "T12" = 0.0001165350087938219
    "G13" = 2.056881333333333
    "H13" = 6.170643999999999
    "I13" = 0.009974175739362347
    "J13" = 0.009974175739362347
    "K13" = 3
    "L13" = 1
    "M13" = 0.4927966666666667
    "N13" = 1.47839
    "O13" = 0.02111942117766284
    "P13" = 0.02111942117766283
    "Q13" = 1.013624264795556
    "R13" = 9.122618383159999
    "S13" = 0.0004706265790900645
    "T13" = 0.0004706265790900644
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output ("Updated {0} cells" -f $updates.Count)